$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, centered) from the existing
# last header cell (G1) onto the new header cell (H1), then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the new "Save" data value for row 2 (plain number, default style).
$ws.Range("H2").Value = 1

$excel.CutCopyMode = 0
